{"js": "// Fix typo: \"Flat Pointer\" -> \"Reference Point\" in the installation instructions\n// paragraph (the app is called \"Reference Point\" everywhere else in the document).\n//\n// We search including the trailing space (\"Flat Pointer \") so the inserted text keeps\n// the space attached to \"Reference Point \" rather than letting it dangle at the start\n// of the next run (\"is a Windows Forms executable...\").\nlet searchResults = context.document.body.search(\"Flat Pointer \", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  for (let i = 0; i < searchResults.items.length; i++) {\n    searchResults.items[i].insertText(\"Reference Point \", Word.InsertLocation.replace);\n  }\n} else {\n  // Fallback in case the trailing space isn't part of the same text run (e.g. a line\n  // break right after \"Pointer\"); just replace the two words themselves.\n  let fallback = context.document.body.search(\"Flat Pointer\", { matchCase: true, matchWholeWord: false });\n  fallback.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < fallback.items.length; i++) {\n    fallback.items[i].insertText(\"Reference Point\", Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fix typo: \"Flat Pointer\" -> \"Reference Point\" in the installation instructions\n# paragraph (Requirements section).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Flat Pointer\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Reference Point\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
